# Add bash entry on "cut" to sheet1 (工作表1) of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$descC = "* cut can accept data passed from pipe or as parameter`n    `$ cat abc.txt | cut -d':' -f1-5`n    `$ cut -d' ' -f1-5 abc.txt`n* cut fix length `n    `$ cut -c1-5 file.txt`n    `$ cut -c10- file.txt`n* cut by delimiter (1-digit) and select fields with option f`n    `$ cut -d':' -f5`n    `$ cut -d':' -f2-6"

$ws.Cells.Item(41, 1).Value = "cut"
$ws.Cells.Item(41, 2).Value = "cut basic"
$ws.Cells.Item(41, 3).Value = $descC

$ws.Rows.Item(41).RowHeight = 141.75
$ws.Cells.Item(41, 3).WrapText = $true

$ws.Range("C42").Select()
